# Updated cryptos list on Sun Mar 26 03:29:17 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for each coin row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "Price" values are numeric-looking strings (e.g. "1.003"), so
# force the cell to Text format before assignment to stop Excel
# auto-coercing them to numbers, then restore the default "Normal" style
# so no stray style index is left referenced on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.644.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.755.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4506"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3556"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07473"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.086"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.25%  "
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.986"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.156"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.750.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001058"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06468"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.762"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.686.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.099"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("E27").Value = "  -1.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.956.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.082"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.084"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09168"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.653"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.497"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02291"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.69"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06028"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("E38").Value = "  -1.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6299"
$ws.Range("D39").Style = "Normal"
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.182"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.86%  "
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.757"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5871"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.940"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06897"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.129"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.08%  "